# Implemented CrossVal for task 3
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 ("Task 3: CrossVal for regression (as in task 1b)") - work started,
# mark status as "in progress" (amber fill) and leave a comment.
$ws.Range("D18").Interior.Color = 49407
$ws.Range("E18").Value = "Started, still a lot of options left to play with"

# Row 19 ("Task 3: Standardize input feature (for regularization)") - also
# started, same "in progress" status color, with a note about what's left.
$ws.Range("D19").Interior.Color = 49407
$ws.Range("E19").Value = "Need to change standardization in CrossVal, right now it's wrong implem"

# Row 9 ("Look for correlation between measurements...") gets a follow-up
# comment about a possible next step.
$ws.Range("E9").Value = "Maybe use feature_selection library to extract best features"

# Reflect the last active selection when the file was saved.
[void]$ws.Range("B13").Select()

Write-Output "done"
